# Apply the cryptos.xlsx data refresh described in the commit
# 'Updated cryptos list on Wed Jun 26 10:51:02 UTC 2024 with GitHub Actions'
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "573.73", "1.00").
# Excel would otherwise silently convert these into real numbers and drop
# formatting such as trailing zeros, so force the cell to Text format first
# so the literal string is preserved exactly, matching the source feed.
$textForceCells = @(
    "D5",
    "D6",
    "D9",
    "D10",
    "D16",
    "D19",
    "D21",
    "D22",
    "D25",
    "D27",
    "D28",
    "D29",
    "D30",
    "D35",
    "D39",
    "D40",
    "D42",
    "D45",
    "D47",
    "D48"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Per-cell values taken from the refreshed coinranking.com feed
$ws.Range("D2").Value = '61.226.01'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").Value = '3.372.06'
$ws.Range("E3").Value = '  -0.02%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '573.73'
$ws.Range("E5").Value = '  +0.42%  '
$ws.Range("D6").Value = '136.47'
$ws.Range("E6").Value = '  -0.76%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = '3.370.51'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.471'
$ws.Range("E9").Value = '  -0.52%  '
$ws.Range("D10").Value = '7.46'
$ws.Range("E10").Value = '  -1.69%  '
$ws.Range("E11").Value = '  +1.15%  '
$ws.Range("E12").Value = '  -0.87%  '
$ws.Range("D13").Value = '3.947.93'
$ws.Range("E13").Value = '  +0.25%  '
$ws.Range("E14").Value = '  +2.43%  '
$ws.Range("E15").Value = '  +1.55%  '
$ws.Range("D16").Value = '26.02'
$ws.Range("E16").Value = '  +3.04%  '
$ws.Range("D17").Value = '3.366.90'
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '61.329.54'
$ws.Range("E18").Value = '  +0.03%  '
$ws.Range("D19").Value = '14.05'
$ws.Range("E19").Value = '  +0.45%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").Value = '9.31'
$ws.Range("E21").Value = '  -1.21%  '
$ws.Range("D22").Value = '375.30'
$ws.Range("E22").Value = '  -1.84%  '
$ws.Range("E23").Value = '  -3.62%  '
$ws.Range("D24").Value = '3.508.36'
$ws.Range("E24").Value = '  +0.09%  '
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.09%  '
$ws.Range("E26").Value = '  +7.21%  '
$ws.Range("D27").Value = '71.43'
$ws.Range("E27").Value = '  +0.82%  '
$ws.Range("D28").Value = '1.70'
$ws.Range("E28").Value = '  +3.36%  '
$ws.Range("B29").Value = 'RenderToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D29").Value = '7.45'
$ws.Range("E29").Value = '  -4.23%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").Value = '1.00'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("E31").Value = '  +1.26%  '
$ws.Range("E32").Value = '  +3.11%  '
$ws.Range("E33").Value = '  +1.45%  '
$ws.Range("E34").Value = '  +0.08%  '
$ws.Range("D35").Value = '23.54'
$ws.Range("E35").Value = '  +0.35%  '
$ws.Range("E36").Value = '  -5.11%  '
$ws.Range("E37").Value = '  -1.91%  '
$ws.Range("D39").Value = '165.43'
$ws.Range("E39").Value = '  +1.22%  '
$ws.Range("D40").Value = '0.0775'
$ws.Range("E40").Value = '  -3.42%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").Value = '0.773'
$ws.Range("E42").Value = '  +1.56%  '
$ws.Range("E43").Value = '  +5.13%  '
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '41.39'
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("E46").Value = '  -0.40%  '
$ws.Range("D47").Value = '24.53'
$ws.Range("E47").Value = '  +5.50%  '
$ws.Range("D48").Value = '6.82'
$ws.Range("E48").Value = '  -2.05%  '
$ws.Range("E49").Value = '  -2.68%  '
$ws.Range("D50").Value = '2.347.77'
$ws.Range("E50").Value = '  +3.39%  '
$ws.Range("E51").Value = '  +0.18%  '
